# contrato_propiedad_definitiva_pagos_varios.docx
# Applies the five wording corrections from the commit:
#   1) "DECLARA "  -> "DECLARAN "  (second DECLARA, the "PROMITENTES" / plural buyers paragraph)
#   2) " y "       -> " Y "        (... DE {{LUGAR_ORIGEN}} y {{LUGAR_ORIGEN_2}}, DE ESTADO CIVIL ...)
#   3) remove the stray period after {{DIRECCION_PROYECTO_LOTE}}
#   4) "EFECTUARÁ " -> "EFECTUARÁN "
#   5) "...{{SEXO_10}}" DEBERÁ CUBRIR..." -> "...DEBERÁN CUBRIR..."

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1
$wdFindContinue = 1
$wdReplaceOne   = 1

# --- 1) Disambiguate the two "DECLARA " occurrences ------------------------
# The first DECLARA belongs to the seller (singular) paragraph and must stay
# untouched; the second belongs to the buyer(s) ("PROMITENTES") paragraph and
# is the one that becomes "DECLARAN ". Anchor the search right after the
# {{CLAUSULA_B}} marker that immediately precedes the second occurrence so
# only it is matched.
$marker = $d.Content
$okMarker = $marker.Find.Execute("{{CLAUSULA_B}}", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $okMarker) { throw "anchor '{{CLAUSULA_B}}' not found" }
$tail = $d.Range($marker.End, $d.Content.End)
$okDeclara = $tail.Find.Execute("DECLARA ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "DECLARAN ", $wdReplaceOne)
if (-not $okDeclara) { throw "'DECLARA ' (buyers' paragraph) not found" }

# --- 2) " y " -> " Y " (only one such lowercase conjunction in the doc) ----
$okY = $d.Content.Find.Execute(" y ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, " Y ", $wdReplaceOne)
if (-not $okY) { throw "' y ' not found" }

# --- 3) Drop the period between the lot address and "EL CUAL" -------------
$okPeriod = $d.Content.Find.Execute("{{DIRECCION_PROYECTO_LOTE}}. EL CUAL", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "{{DIRECCION_PROYECTO_LOTE}} EL CUAL", $wdReplaceOne)
if (-not $okPeriod) { throw "'{{DIRECCION_PROYECTO_LOTE}}. EL CUAL' not found" }

# --- 4) "EFECTUARÁ " -> "EFECTUARÁN " --------------------------------------
$okEfectuara = $d.Content.Find.Execute("EFECTUARÁ ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "EFECTUARÁN ", $wdReplaceOne)
if (-not $okEfectuara) { throw "'EFECTUARÁ ' not found" }

# --- 5) "...DEBERÁ CUBRIR LA PENA CONVENCIONAL ESTABLECIDA..." -> "...DEBERÁN CUBRIR..." ---
$okDebera = $d.Content.Find.Execute("DEBERÁ CUBRIR LA PENA CONVENCIONAL ESTABLECIDA", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "DEBERÁN CUBRIR LA PENA CONVENCIONAL ESTABLECIDA", $wdReplaceOne)
if (-not $okDebera) { throw "'DEBERÁ CUBRIR LA PENA CONVENCIONAL ESTABLECIDA' not found" }

Write-Output "OK: 5/5 replacements applied"
